# This workbook holds weekly price-report rows (2-15) for the same market/product
# that simply got reshuffled to a different week ordering. Column D (Fecha) and
# columns K..T (Variedad..Kg/unidad) move together as a unit; columns A..C and
# E..J are identical for every row, so only D and K:T need to be relocated.
#
# Mapping of resulting row -> source row (1-based worksheet rows):
#   2 <- 3, 3 <- 4, 4 <- 13, 5 <- 10, 6 <- 14, 7 <- 11, 8 <- 12,
#   9 <- 7, 10 <- 15, 11 <- 6, 12 <- 8, 13 <- 9, 14 <- 2, 15 <- 5
# (rows 1 and 16 stay as-is)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 3
    3  = 4
    4  = 13
    5  = 10
    6  = 14
    7  = 11
    8  = 12
    9  = 7
    10 = 15
    11 = 6
    12 = 8
    13 = 9
    14 = 2
    15 = 5
}

# Columns that move as a block: D (4) and K..T (11..20)
$cols = @(4,11,12,13,14,15,16,17,18,19,20)

# Snapshot every source cell's value before any writes happen, so that later
# writes never clobber a value that is still needed as a source for another row.
$snapshot = @{}
for ($r = 2; $r -le 15; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $snapshot["$srcRow-$c"]
    }
}
